$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label / data cell values ---
$ws.Range("B1").Value = "Element"
$ws.Range("C1").Value = "Biorefinery"
$ws.Range("B2").Value = "Metric"
$ws.Range("C2").Value = "MFPP [USD/ton]"
$ws.Range("D2").Value = "Biodiesel production [MMGal/yr]"
$ws.Range("E2").Value = "Ethanol production [MMGal/yr]"
$ws.Range("F2").Value = "Electricity production [MMWhr/yr]"
$ws.Range("G2").Value = "Natural gas consumption [MMcf/yr]"
$ws.Range("H2").Value = "Productivity [MMGGE/yr]"
$ws.Range("I2").Value = "TCI [10^6*USD]"
$ws.Range("J2").Value = "Feedstock consumption [ton/yr]"
$ws.Range("K2").Value = "Heat exchanger network error [%]"
$ws.Range("A3").Value = "Element"
$ws.Range("B3").Value = "Parameter"
$ws.Range("A4").Value = "Stream-lipidcane"
$ws.Range("B4").Value = "Lipid content [dry wt. %]"
$ws.Range("C4").Value = 0.2825275631890798
$ws.Range("D4").Value = 0.9586439660991526
$ws.Range("E4").Value = -0.7866451661291534
$ws.Range("G4").Value = -0.3830765769144229
$ws.Range("H4").Value = 0.5490497262431562
$ws.Range("I4").Value = 0.3512847821195531
$ws.Range("J4").Value = -0.07455936398409961
$ws.Range("K4").Value = -0.252172804320108
$ws.Range("B5").Value = "Lipid retention [%]"
$ws.Range("C5").Value = 0.02469961749043727
$ws.Range("D5").Value = 0.09077476936923425
$ws.Range("E5").Value = 0.02245256131403285
$ws.Range("G5").Value = 0.007884197104927625
$ws.Range("H5").Value = 0.1661666541663542
$ws.Range("I5").Value = 0.13750393759844
$ws.Range("J5").Value = 0.1357848946223656
$ws.Range("K5").Value = -0.1363849096227406
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"
$ws.Range("C6").Value = 0.1011790294757369
$ws.Range("D6").Value = 0.1886012150303758
$ws.Range("E6").Value = -0.0254166354158854
$ws.Range("G6").Value = 0.5542203555088878
$ws.Range("H6").Value = -0.03332333308332708
$ws.Range("I6").Value = -0.02942173554338859
$ws.Range("J6").Value = -0.06317257931448288
$ws.Range("K6").Value = -0.07991449786244657
$ws.Range("B7").Value = "Capacity [ton/hr]"
$ws.Range("C7").Value = 0.07266481662041552
$ws.Range("D7").Value = 0.08541513537838447
$ws.Range("E7").Value = 0.5765364134103353
$ws.Range("G7").Value = 0.3828515712892822
$ws.Range("H7").Value = 0.7021705542638568
$ws.Range("I7").Value = 0.8934523363084078
$ws.Range("J7").Value = 0.8632430810770271
$ws.Range("K7").Value = 0.3285187129678242
$ws.Range("A8").Value = "Stream-ethanol"
$ws.Range("B8").Value = "Price [USD/gal]"
$ws.Range("C8").Value = 0.802236555913898
$ws.Range("D8").Value = 0.04826070651766294
$ws.Range("E8").Value = 0.0003645091127278182
$ws.Range("G8").Value = -0.011716792919823
$ws.Range("H8").Value = 0.07834995874896875
$ws.Range("I8").Value = 0.06283957098927474
$ws.Range("J8").Value = 0.04300307507687693
$ws.Range("K8").Value = -0.07271731793294833
$ws.Range("A9").Value = "Stream-biodiesel"
$ws.Range("B9").Value = "Price [USD/gal]"
$ws.Range("C9").Value = 0.3791644791119778
$ws.Range("D9").Value = -0.09179029475736895
$ws.Range("E9").Value = 0.0381969549238731
$ws.Range("G9").Value = 0.003286582164554114
$ws.Range("H9").Value = -0.05133128328208206
$ws.Range("I9").Value = -0.003903097577439436
$ws.Range("J9").Value = -0.005701642541063528
$ws.Range("K9").Value = 0.01346583664591615
$ws.Range("A10").Value = "Stream-natural gas"
$ws.Range("B10").Value = "Price [USD/cf]"
$ws.Range("C10").Value = 0.02707267681692043
$ws.Range("D10").Value = 0.01770794269856747
$ws.Range("E10").Value = -0.04158853971349284
$ws.Range("G10").Value = -0.05002775069376735
$ws.Range("H10").Value = 0.00172654316357909
$ws.Range("I10").Value = -0.02691667291682293
$ws.Range("J10").Value = -0.01775744393609841
$ws.Range("K10").Value = 0.0716327908197705
$ws.Range("A11").Value = "biorefinery"
$ws.Range("B11").Value = "Electricity price [USD/kWh]"
$ws.Range("C11").Value = -0.00666316657916448
$ws.Range("D11").Value = 0.08472511812795321
$ws.Range("E11").Value = 0.03169729243231081
$ws.Range("G11").Value = 0.05180529513237832
$ws.Range("H11").Value = 0.1112832820820521
$ws.Range("I11").Value = 0.1310837770944274
$ws.Range("J11").Value = 0.1039225980649516
$ws.Range("K11").Value = 0.02524113102827571
$ws.Range("B12").Value = "Operating days [day/yr]"
$ws.Range("C12").Value = 0.04452861321533038
$ws.Range("D12").Value = 0.07727293182329559
$ws.Range("E12").Value = 0.2624705617640442
$ws.Range("G12").Value = 0.5191224780619517
$ws.Range("H12").Value = 0.224537613440336
$ws.Range("I12").Value = -0.08857871446786171
$ws.Range("J12").Value = 0.401225530638266
$ws.Range("K12").Value = 0.1231050776269407
$ws.Range("B13").Value = "IRR [%]"
$ws.Range("C13").Value = -0.1908212705317633
$ws.Range("D13").Value = 0.004089102227555689
$ws.Range("E13").Value = -0.01011175279381985
$ws.Range("G13").Value = -0.0229940748518713
$ws.Range("H13").Value = 0.04297307432685817
$ws.Range("I13").Value = 0.06355058876471913
$ws.Range("J13").Value = 0.01811745293632341
$ws.Range("K13").Value = -0.08416260406510165

# --- Match header/styled formatting for the new column K cells (row1 blank, row2 label) ---
$ws.Range("K2").Font.Bold = $true
$ws.Range("K2").HorizontalAlignment = -4108
$ws.Range("K2").VerticalAlignment = -4160
$ws.Range("K2").Borders.LineStyle = 1

# --- Extend merged header range C1:J1 -> C1:K1, keep uniform header styling ---
$ws.Range("C1:K1").Merge()
$hdr = $ws.Range("C1:K1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
